$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: Price (D) and Volume 1h change (E) columns.
# Values are stored as text (matching the source data which is inline
# string/text, not numeric) so number-formatted strings like
# "28.488.32" or "1.000" are not reinterpreted as numbers/dates.
$updates = @(
    @{ Row = 2; D = "28.488.32"; E = "  -0.81%  " },
    @{ Row = 3; D = "1.830.10"; E = "  +1.23%  " },
    @{ Row = 4; D = "1.002"; E = "  -0.02%  " },
    @{ Row = 5; D = "329.79"; E = "  +0.53%  " },
    @{ Row = 6; D = "1.000"; E = "  +0.08%  " },
    @{ Row = 7; D = "0.4524"; E = "  +3.23%  " },
    @{ Row = 8; D = "0.3804"; E = "  +0.90%  " },
    @{ Row = 9; D = "44.89"; E = "  -0.81%  " },
    @{ Row = 10; D = "0.07811"; E = "  +1.61%  " },
    @{ Row = 11; D = "1.144"; E = "  +0.36%  " },
    @{ Row = 12; D = "22.37"; E = "  -1.38%  " },
    @{ Row = 13; D = "1.001"; E = "  -0.05%  " },
    @{ Row = 14; D = "6.399"; E = "  +2.07%  " },
    @{ Row = 15; D = "7.552"; E = "  +0.50%  " },
    @{ Row = 16; D = "1.837.45"; E = "  +1.55%  " },
    @{ Row = 17; D = "94.17"; E = "  +16.13%  " },
    @{ Row = 18; E = "  -0.33%  " },
    @{ Row = 19; D = "0.06389"; E = "  -4.92%  " },
    @{ Row = 20; D = "0.9999"; E = "  +0.01%  " },
    @{ Row = 21; D = "17.61"; E = "  -0.21%  " },
    @{ Row = 22; D = "6.402"; E = "  +1.87%  " },
    @{ Row = 23; D = "0.5423"; E = "  -1.10%  " },
    @{ Row = 24; D = "28.540.58"; E = "  -0.59%  " },
    @{ Row = 25; D = "11.83"; E = "  +0.48%  " },
    @{ Row = 26; D = "2.298"; E = "  -6.22%  " },
    @{ Row = 27; D = "20.97"; E = "  +1.95%  " },
    @{ Row = 28; D = "153.70"; E = "  -0.68%  " },
    @{ Row = 29; D = "2.373"; E = "  +0.49%  " },
    @{ Row = 30; D = "2.043.77"; E = "  +1.37%  " },
    @{ Row = 31; D = "129.69"; E = "  -1.01%  " },
    @{ Row = 32; D = "1.210"; E = "  -7.50%  " },
    @{ Row = 33; D = "5.908"; E = "  +1.37%  " },
    @{ Row = 34; D = "0.09338"; E = "  +1.47%  " },
    @{ Row = 35; D = "3.668"; E = "  -7.55%  " },
    @{ Row = 36; D = "12.94"; E = "  +6.06%  " },
    @{ Row = 37; D = "0.02366"; E = "  +2.29%  " },
    @{ Row = 38; D = "0.2210"; E = "  -0.70%  " },
    @{ Row = 39; D = "0.6698"; E = "  +1.34%  " },
    @{ Row = 40; D = "0.06304"; E = "  -0.27%  " },
    @{ Row = 41; D = "5.233"; E = "  +0.50%  " },
    @{ Row = 42; D = "8.191"; E = "  +1.55%  " },
    @{ Row = 43; D = "1.201"; E = "  -0.36%  " },
    @{ Row = 44; B = "WEMIXTOKEN"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D = "1.411"; E = "  -1.87%  " },
    @{ Row = 45; D = "0.9998"; E = "  +0.10%  " },
    @{ Row = 46; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "14.00"; E = "  -0.16%  " },
    @{ Row = 47; D = "0.6174"; E = "  +1.57%  " },
    @{ Row = 48; D = "3.779"; E = "  -0.40%  " },
    @{ Row = 49; D = "128.24"; E = "  +0.31%  " },
    @{ Row = 50; D = "2.058"; E = "  +1.57%  " },
    @{ Row = 51; E = "  -0.76%  " }
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    if ($item.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    if ($item.ContainsKey("D")) {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
        $ws.Cells.Item($r, 4).Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).NumberFormat = "@"
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}

